# Plan tests acceptation - add the "confirmation" row content (row 8)
# and update the saved view selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the four cells of row 8 with their text values.
$ws.Range("B8").Value = "Un formulaire permet à l'utilisateur de valider la commande !"
$ws.Range("C8").Value = "La validité des données rentrées dans les champs du formulaire est testée en direct, et une seconde fois lors du clic sur le bouton ""Commander !"""
$ws.Range("D8").Value = "Si les champs du formulaire sont tous valides, la commande peut être envoyée après un message de confirmation comportant un récapitulatif des données du formulaire, du nombre d'articles commandés, ainsi que le prix total.  "
$ws.Range("E8").Value = "OK / Si la connexion avec l'API ne peut être établie, l'utilisateur en sera informé par un message d'alerte. Si le panier est vide, l'utilisateur sera informé par une alerte de l'impossibilité de passer la commande."

# Row 8 grows taller to fit the new text (same height pattern as other filled rows).
$ws.Rows.Item(8).RowHeight = 129.6

# Update the view: scrolled to column C, selection now on E8.
$ws.Range("E8").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 7
